$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-13) from serial date 45204 to 45207
$ws.Range("C2:C13").Value = 45207
